# Add two new time-entries (2014-03-31) before the summary rows.
# This mirrors inserting two data rows just above the existing "sum [min]"
# row, pushing the summary block (sum [min] / sum [h] / sum [working weeks])
# down by two rows and widening the SUM() range that feeds it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh rows right before the old row 104 (the blank template row
# that precedes the summary block). The summary rows below shift from
# 105-107 to 107-109 and their formulas/ranges adjust automatically.
$ws.Rows("104:105").Insert()

# Row 104: 2014-03-31, 10:15 -> 14:00
$ws.Range("A104").Value = 2014
$ws.Range("B104").Value = 3
$ws.Range("C104").Value = 31
$ws.Range("D104").Value = 0.42708333333333331
$ws.Range("E104").Value = 0.58333333333333337
$ws.Range("F104").Formula = "=(E104-D104)*24*60"
$ws.Range("G104").Formula = "=F104/60"

# Row 105: 2014-03-31, 14:15 -> 16:00
$ws.Range("A105").Value = 2014
$ws.Range("B105").Value = 3
$ws.Range("C105").Value = 31
$ws.Range("D105").Value = 0.59375
$ws.Range("E105").Value = 0.66666666666666663
$ws.Range("F105").Formula = "=(E105-D105)*24*60"
$ws.Range("G105").Formula = "=F105/60"

# Move the selection to A106 (the now-blank spacer row above the summary).
$null = $ws.Range("A106").Select()
